$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-29 10:49:26"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-29 10:49:22"
$wsZhCn.Range("K2").Value = "2016-08-29 10:49:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-29 10:49:26"
$wsDeDe.Range("K2").Value = "2016-08-29 10:49:46"
